# Update the four "Munroe Score" ovals on the single slide: the score
# values are refreshed to more precise numbers and the font size is
# reduced from 40pt to 28pt so the longer (decimal) strings still fit
# neatly inside the ovals.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Oval 6" was originally typed as two separate runs ("6" + "0"). Drop
# the leading run first so retyping the value collapses everything back
# down to a single run - matching how PowerPoint merges runs when the
# whole field is retyped at once - while keeping the trailing run (and
# its smtClean marker) alive.
$oval6 = $s.Shapes.Item("Oval 6")
$oval6.TextFrame.TextRange.Characters(1, 1).Delete()

$newScore = 28

$oval = $s.Shapes.Item("Oval 6")
$tr = $oval.TextFrame.TextRange
$tr.Text = "55.1"
$tr.Font.Size = $newScore

$oval = $s.Shapes.Item("Oval 9")
$tr = $oval.TextFrame.TextRange
$tr.Text = "59.4"
$tr.Font.Size = $newScore

$oval = $s.Shapes.Item("Oval 10")
$tr = $oval.TextFrame.TextRange
$tr.Text = "67.5"
$tr.Font.Size = $newScore

$oval = $s.Shapes.Item("Oval 11")
$tr = $oval.TextFrame.TextRange
$tr.Text = "61.4"
$tr.Font.Size = $newScore
